$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows where "Recorded By" (column G) was "System, backup@backdoor.com, system"
# and becomes "system, System, backup@backdoor.com"
$rowsTypeA = @(2, 28, 54)
foreach ($r in $rowsTypeA) {
    $ws.Cells.Item($r, 7).Value = "system, System, backup@backdoor.com"
}

# Rows where "Recorded By" (column G) was "dnasr281@gmail.com, System"
# and becomes "System, dnasr281@gmail.com"
$rowsTypeB = @(10, 11, 12, 13, 14, 15, 17, 18, 19, 20, 21, 22, 24, 26, 36, 37, 38, 39, 40, 41, 43, 44, 45, 46, 47, 48, 50, 52, 62, 63, 64, 65, 66, 67, 69, 70, 71, 72, 73, 74, 76, 78, 83, 84, 85, 86, 90, 92, 93, 94, 96, 99, 101, 109, 110, 111, 112, 116, 118, 119, 120, 122, 125, 127, 135, 136, 137, 138, 142, 144, 145, 146, 148, 151, 153)
foreach ($r in $rowsTypeB) {
    $ws.Cells.Item($r, 7).Value = "System, dnasr281@gmail.com"
}
